# Append the new purchase row (row 48) recorded after running on 2025-08-08.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (matching the existing rows'
# "MM/DD/YYYY" inline-string convention) -- the leading apostrophe forces
# Excel to store it as text instead of auto-converting it to a date serial.
$ws.Cells.Item(48, 1).Value = "'08/08/2025"
$ws.Cells.Item(48, 2).Value = 552.2389999999941
$ws.Cells.Item(48, 3).Value = 0.09054050872901141
$ws.Cells.Item(48, 4).Value = 50
